$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 80
$ws.Range("F6").Value = 728
$ws.Range("F7").Value = 186
$ws.Range("F8").Value = 228
$ws.Range("F11").Value = 43
$ws.Range("F12").Value = 579
$ws.Range("F13").Value = 488
$ws.Range("F16").Value = 142
$ws.Range("F17").Value = 798
$ws.Range("F18").Value = 2585
$ws.Range("F23").Value = 195
$ws.Range("F25").Value = 142
$ws.Range("F27").Value = 949
$ws.Range("F29").Value = 184
$ws.Range("F30").Value = 4
$ws.Range("F34").Value = 266

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 988
$ws.Range("F5").Value = 988
$ws.Range("F10").Value = 297
$ws.Range("F14").Value = 545
$ws.Range("F17").Value = 960
$ws.Range("F24").Value = 271
$ws.Range("F25").Value = 238
$ws.Range("F31").Value = 76

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1763
$ws.Range("F5").Value = 2368
$ws.Range("F6").Value = 960
$ws.Range("F9").Value = 1203

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1763
$ws.Range("F4").Value = 2368
$ws.Range("F8").Value = 960
$ws.Range("F9").Value = 1203
$ws.Range("F13").Value = 80
$ws.Range("F14").Value = 728
$ws.Range("F15").Value = 186
$ws.Range("F17").Value = 228
$ws.Range("F19").Value = 43
$ws.Range("F20").Value = 579
$ws.Range("F21").Value = 988
$ws.Range("F22").Value = 488
$ws.Range("F25").Value = 142
$ws.Range("F26").Value = 798
$ws.Range("F27").Value = 2585
$ws.Range("F31").Value = 195
$ws.Range("F32").Value = 142
$ws.Range("F34").Value = 949
$ws.Range("F35").Value = 545
$ws.Range("F38").Value = 184
$ws.Range("F43").Value = 271
$ws.Range("F44").Value = 271
$ws.Range("F45").Value = 238
$ws.Range("F50").Value = 266

